$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix wrong subject codes (bug fix mentioned in commit message)
$ws.Range("C13").Value = "QP2"
$ws.Range("C14").Value = "Tin5"

# Update the active cell selection saved with the sheet
$ws.Range("C14").Select()
